$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3-9 down to 4-10.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with this week's data.
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "Vega Monumental Concepción"
$ws.Range("C3").Value = "Bíobío"
$ws.Range("D3").Value = 44496
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 100112022
$ws.Range("G3").Value = "Arveja Verde"
$ws.Range("H3").Value = "Perfection"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 14520
$ws.Range("N3").Value = "$/malla 25 kilos"
$ws.Range("O3").Value = "Provincia de Huasco"
$ws.Range("P3").Value = 581
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
